# Natmi following Dr Hou advice
#
# The LR-pair sheet (Dll1 -> Notch3) was re-derived: the 3 original rows
# (Sending cluster = sCs only) are replaced with updated numbers, and the
# pipeline now also reports the symmetric case where the sending cluster
# is ECs, for a total of 6 data rows (2-7) under the header row (row 1).
#
# Columns (A..T):
#   A Sending cluster            K Receptor-expressing cells
#   B Ligand symbol              L Receptor detection rate
#   C Receptor symbol            M Receptor average expression value
#   D Target cluster             N Receptor total expression value
#   E Ligand-expressing cells    O Receptor derived specificity (avg)
#   F Ligand detection rate      P Receptor derived specificity (total)
#   G Ligand average expr value  Q Edge average expression weight
#   H Ligand total expr value    R Edge total expression weight
#   I Ligand derived spec (avg)  S Edge average expression derived specificity
#   J Ligand derived spec (tot)  T Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @{
    2 = @("ECs","Dll1","Notch3","ECs",2,0.6666666666666666,10.34874566666667,31.046237,0.614862320492409,0.6148623204924089,3,1,7.457778999999999,22.373337,0.08029647035915141,0.0802964703591514,77.17865810920767,694.607922982869,0.04937127409237777,0.04937127409237776)
    3 = @("ECs","Dll1","Notch3","FAPs",2,0.6666666666666666,10.34874566666667,31.046237,0.614862320492409,0.6148623204924089,3,1,2.732509,8.197526999999999,0.02942039820764526,0.02942039820764525,28.27804067287767,254.502366055899,0.01808949431176348,0.01808949431176347)
    4 = @("ECs","Dll1","Notch3","sCs",2,0.6666666666666666,10.34874566666667,31.046237,0.614862320492409,0.6148623204924089,3,1,82.68775466666666,248.063264,0.8902831314332034,0.8902831314332033,855.7145427930632,7701.430885137568,0.5474015520882678,0.5474015520882676)
    5 = @("sCs","Dll1","Notch3","ECs",3,1,6.482251000000001,19.446753,0.3851376795075911,0.385137679507591,3,1,7.457778999999999,22.373337,0.08029647035915141,0.0802964703591514,48.343195380529,435.088758424761,0.03092519626677364,0.03092519626677363)
    6 = @("sCs","Dll1","Notch3","FAPs",3,1,6.482251000000001,19.446753,0.3851376795075911,0.385137679507591,3,1,2.732509,8.197526999999999,0.02942039820764526,0.02942039820764525,17.712809197759,159.415282779831,0.01133090389588179,0.01133090389588178)
    7 = @("sCs","Dll1","Notch3","sCs",3,1,6.482251000000001,19.446753,0.3851376795075911,0.385137679507591,3,1,82.68775466666666,248.063264,0.8902831314332034,0.8902831314332033,536.0027803757547,4824.025023381792,0.3428815793449357,0.3428815793449356)
}

foreach ($r in 2..7) {
    $rowVals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + "$r").Value = $rowVals[$i]
    }
}
